# Fruta / hortaliza, semanal
# The weekly refresh swaps the data that used to live in rows 2-3 (week of
# 44505) with the data that used to live in rows 4-5 (week of 44902), so the
# most recent week now appears first. Columns A, B, C, E, F, G, H, I, J, R
# and T are identical across all four data rows and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes what used to be row 4 ---
$ws.Range("D2").Value = 44902
$ws.Range("K2").Value = "Golden Nugget"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/caja 10 kilos"
$ws.Range("S2").Value = 1500

# --- Row 3 becomes what used to be row 5 ---
$ws.Range("D3").Value = 44902
$ws.Range("K3").Value = "Golden Nugget"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("S3").Value = 1300

# --- Row 4 becomes what used to be row 2 ---
$ws.Range("D4").Value = 44505
$ws.Range("K4").Value = "Californiana(o)"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1500

# --- Row 5 becomes what used to be row 3 ---
$ws.Range("D5").Value = 44505
$ws.Range("K5").Value = "Golden Nugget"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("S5").Value = 1500
